# Auto-generated Excel COM-interop script to apply Bahamut_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 196.4
$ws.Range("I33").Value = 289.5
$ws.Range("J33").Value = 134.33333
$ws.Range("K33").Value = 289.5
$ws.Range("L33").Value = 134.33333
$ws.Range("M33").Value = -60.5
$ws.Range("N33").Value = -592.3333299999999
$ws.Range("H64").Value = 3900
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9752
$ws.Range("H67").Value = 3900
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9142
$ws.Range("H74").Value = 12734147
$ws.Range("I74").Value = 12734147
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 12734147
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -12733211
$ws.Range("H76").Value = 11113944
$ws.Range("I76").Value = 33335834
$ws.Range("K76").Value = 33335834
$ws.Range("M76").Value = -33335519
$ws.Range("H77").Value = 12734147
$ws.Range("I77").Value = 12734147
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 63670735
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -63666055
$ws.Range("H79").Value = 11113944
$ws.Range("I79").Value = 33335834
$ws.Range("K79").Value = 33335834
$ws.Range("M79").Value = -33334742
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2300.2942
$ws.Range("I2").Value = 2140.3333
$ws.Range("J2").Value = 3500
$ws.Range("K2").Value = 2140.3333
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = -2027.3333
$ws.Range("N2").Value = -3726
$ws.Range("H45").Value = 2636.5715
$ws.Range("I45").Value = 1355.6364
$ws.Range("K45").Value = 1355.6364
$ws.Range("M45").Value = -978.6364000000001
$ws.Range("H61").Value = 1010.2778
$ws.Range("I61").Value = 742.6429000000001
$ws.Range("K61").Value = 742.6429000000001
$ws.Range("M61").Value = -530.6429000000001
$ws.Range("H88").Value = 2284
$ws.Range("I88").Value = 1821.8
$ws.Range("J88").Value = 3439.5
$ws.Range("K88").Value = 1821.8
$ws.Range("L88").Value = 3439.5
$ws.Range("M88").Value = -1415.8
$ws.Range("N88").Value = -4251.5
$ws.Range("H91").Value = 2284
$ws.Range("I91").Value = 1821.8
$ws.Range("J91").Value = 3439.5
$ws.Range("K91").Value = 1821.8
$ws.Range("L91").Value = 3439.5
$ws.Range("M91").Value = -417.8
$ws.Range("N91").Value = -6247.5
$ws.Range("H97").Value = 709.75
$ws.Range("I97").Value = 425.875
$ws.Range("J97").Value = 1277.5
$ws.Range("K97").Value = 425.875
$ws.Range("L97").Value = 1277.5
$ws.Range("M97").Value = 70.125
$ws.Range("N97").Value = -2269.5
$ws.Range("H116").Value = 2300.2942
$ws.Range("I116").Value = 2140.3333
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2140.3333
$ws.Range("L116").Value = 3500
$ws.Range("M116").Value = 153.6667000000002
$ws.Range("N116").Value = -8088
$ws.Range("H136").Value = 1010.2778
$ws.Range("I136").Value = 742.6429000000001
$ws.Range("K136").Value = 2227.9287
$ws.Range("M136").Value = 322.0712999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2300.2942
$ws.Range("I3").Value = 2140.3333
$ws.Range("J3").Value = 3500
$ws.Range("K3").Value = 2140.3333
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = -2026.3333
$ws.Range("N3").Value = -3728
$ws.Range("H94").Value = 929.2222
$ws.Range("I94").Value = 948.4
$ws.Range("J94").Value = 833.3333
$ws.Range("K94").Value = 948.4
$ws.Range("L94").Value = 833.3333
$ws.Range("M94").Value = -497.4
$ws.Range("N94").Value = -1735.3333
$ws.Range("H105").Value = 6078.5713
$ws.Range("I105").Value = 5666.6665
$ws.Range("J105").Value = 6820
$ws.Range("K105").Value = 5666.6665
$ws.Range("L105").Value = 6820
$ws.Range("M105").Value = -3919.6665
$ws.Range("N105").Value = -10314
$ws.Range("H107").Value = 6934.5356
$ws.Range("I107").Value = 2128.9443
$ws.Range("K107").Value = 2128.9443
$ws.Range("M107").Value = -208.9443000000001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 24330.4
$ws.Range("J131").Value = 24330.4
$ws.Range("L131").Value = 24330.4
$ws.Range("N131").Value = -34410.4
$ws.Range("H132").Value = 3374.8
$ws.Range("I132").Value = 2190
$ws.Range("J132").Value = 4559.6
$ws.Range("K132").Value = 6570
$ws.Range("L132").Value = 13678.8
$ws.Range("M132").Value = -4040
$ws.Range("N132").Value = -18738.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1155.5186
$ws.Range("I122").Value = 650
$ws.Range("J122").Value = 1195.96
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 10763.64
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -15663.64
$ws.Range("H141").Value = 4141.9414
$ws.Range("I141").Value = 1705.25
$ws.Range("K141").Value = 5115.75
$ws.Range("M141").Value = 64.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3958.4692
$ws.Range("I70").Value = 3971.8413
$ws.Range("J70").Value = 3911.6667
$ws.Range("K70").Value = 3971.8413
$ws.Range("L70").Value = 3911.6667
$ws.Range("M70").Value = -3701.8413
$ws.Range("N70").Value = -4451.6667
$ws.Range("H73").Value = 3958.4692
$ws.Range("I73").Value = 3971.8413
$ws.Range("J73").Value = 3911.6667
$ws.Range("K73").Value = 3971.8413
$ws.Range("L73").Value = 3911.6667
$ws.Range("M73").Value = -3035.8413
$ws.Range("N73").Value = -5783.6667
$ws.Range("H80").Value = 3964.5
$ws.Range("I80").Value = 3963.6365
$ws.Range("J80").Value = 3965.5557
$ws.Range("K80").Value = 3963.6365
$ws.Range("L80").Value = 3965.5557
$ws.Range("M80").Value = -2965.6365
$ws.Range("N80").Value = -5961.5557
$ws.Range("H83").Value = 3964.5
$ws.Range("I83").Value = 3963.6365
$ws.Range("J83").Value = 3965.5557
$ws.Range("K83").Value = 19818.1825
$ws.Range("L83").Value = 19827.7785
$ws.Range("M83").Value = -14826.1825
$ws.Range("N83").Value = -29811.7785
$ws.Range("H97").Value = 1113.2963
$ws.Range("I97").Value = 1264.238
$ws.Range("J97").Value = 585
$ws.Range("K97").Value = 1264.238
$ws.Range("L97").Value = 585
$ws.Range("M97").Value = -768.2380000000001
$ws.Range("N97").Value = -1577
$ws.Range("H132").Value = 6137.125
$ws.Range("I132").Value = 11250
$ws.Range("J132").Value = 4432.8335
$ws.Range("K132").Value = 33750
$ws.Range("L132").Value = 13298.5005
$ws.Range("M132").Value = -31220
$ws.Range("N132").Value = -18358.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 99.666664
$ws.Range("I32").Value = 99.666664
$ws.Range("K32").Value = 99.666664
$ws.Range("M32").Value = 217.333336
$ws.Range("H40").Value = 10101010
$ws.Range("I40").Value = 10101010
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 10101010
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -10100874
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("N40").ClearContents()
